$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert 4 new rows for the additional Section 3 requirements (3.6-3.9).
#    This pushes the existing Section 4 block (row 18-19) down to rows 22-23.
# ---------------------------------------------------------------------------
$ws.Rows("17:20").Insert()

# ---------------------------------------------------------------------------
# 2. Text updates, applied in the same order the requirements were revised
# ---------------------------------------------------------------------------
# Section 4 / Platform - existing requirement reworded
$ws.Range("B23").Value = "User actions shall be performed in a graphical user interface"

# 3.2 requirement text rewritten
$ws.Range("B13").Value = "The system shall calculate the shortest route between all locations, starting at the chosen start waypoint and terminating at the chosen end waypoint"

# 1.2.3 requirement text - capitalize "Admins"
$ws.Range("B6").Value = "Admins shall be able to remove sites"

# Section 2 sample requirement - "locations" -> '"waypoints"'
$ws.Range("B9").Value = 'The user shall be able to identify start and end "waypoints"'

# ---------------------------------------------------------------------------
# 3. Fill in the newly inserted Section 3 rows (3.6 - 3.9)
# ---------------------------------------------------------------------------
$ws.Range("A17").Value = 3.6
$ws.Range("B17").Value = "The system shall include all US national parks as waypoints"

$ws.Range("A18").Value = 3.7
$ws.Range("B18").Value = "The system shall include Machu Picchu as a waypoint"

$ws.Range("A19").Value = 3.8
$ws.Range("B19").Value = "The system shall calculate the final route distance"

$ws.Range("A20").Value = 3.9
$ws.Range("B20").Value = "The system shall include a facts display for each waypoint"

# ---------------------------------------------------------------------------
# 4. New Section 4 requirement (4.2)
# ---------------------------------------------------------------------------
$ws.Range("A24").Value = 4.2
$ws.Range("B24").Value = "The GUI shall be blue and green"

# The Section 4 sample answers (B23/B24) are regular (non-bold) text
$ws.Range("B23").Font.Bold = $false
$ws.Range("B24").Font.Bold = $false

# ---------------------------------------------------------------------------
# 5. Selection state (last edited cell)
# ---------------------------------------------------------------------------
$ws.Range("B24").Select()
